$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new columns before column D (shift existing D:K to F:M).
$ws.Columns("D:E").Insert(-4161)

# 2. Copy number formats from the (now-shifted) old D:E columns, which now
#    live at F:G, onto the freshly inserted D:E columns so the new cells
#    pick up the same date / number styling instead of the generic default.
#    (Restrict to the used row range so we don't balloon the sheet to 1M+ rows.)
$ws.Range("F1:G102").Copy()
$ws.Range("D1:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newDE = @{
    7 = @(43465, 43373)
    8 = @(180800, 182200)
    9 = @(128400, 129800)
    10 = @(52400, 52400)
    11 = @($null, $null)
    12 = @(900, 700)
    13 = @(0, 0)
    14 = @(1700, 3700)
    15 = @(3100, 2900)
    16 = @($null, $null)
    17 = @(178300, 179200)
    18 = @(2500, 3000)
    19 = @($null, $null)
    20 = @(0, 0)
    21 = @(11200, 11500)
    22 = @(2400, 1900)
    23 = @(100, 1100)
    24 = @(0, 2100)
    25 = @(0, 0)
    26 = @(200, -1000)
    27 = @(200, -1000)
    28 = @(0, 0)
    29 = @(-1200, "NA")
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(0, 0)
    33 = @(-1100, -1000)
    34 = @(0, 0)
    35 = @(-1100, -1000)
    38 = @(43465, 43373)
    39 = @($null, $null)
    40 = @($null, $null)
    41 = @(25500, 17100)
    42 = @(0, 0)
    43 = @(148300, 155600)
    44 = @(13100, 11100)
    45 = @(15900, 15600)
    46 = @(202800, 199400)
    47 = @(0, 0)
    48 = @(93900, 86400)
    49 = @(390700, 256100)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(6700, 6200)
    53 = @(0, 0)
    54 = @(694000, 548100)
    55 = @($null, $null)
    56 = @($null, $null)
    57 = @(13900, 12900)
    58 = @(10800, 7300)
    59 = @(75900, 75000)
    60 = @(100500, 95200)
    61 = @(292900, 156400)
    62 = @(29600, 18300)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(423100, 270100)
    67 = @($null, $null)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(71600, 72600)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(270900, 278100)
    77 = @(0, 0)
    80 = @(43465, 43373)
    81 = @(-1100, -1000)
    82 = @($null, $null)
    83 = @(8700, 8500)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(17500, 4100)
    90 = @($null, $null)
    91 = @(-5200, -4400)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(-145600, 500)
    95 = @($null, $null)
    96 = @(0, 0)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(137900, -4600)
    101 = @(-1300, -400)
    102 = @(8500, -500)
}

$rowOverrides = @{
    89 = @(17500, 4100, 14300, 5800, 20600, 12300, 9600, 13400, 4300, 8600)
    91 = @(-5200, -4400, -5800, -5200, -4900, -4600, -6400, -3800, -2300, -3400)
}

# 3. Populate the two new columns (D and E) with the new quarter data for
#    every row of the report.
foreach ($r in $newDE.Keys) {
    $vals = $newDE[$r]
    $ws.Cells.Item([int]$r, 4).Value = $vals[0]
    $ws.Cells.Item([int]$r, 5).Value = $vals[1]
}

# 4. A couple of rows (89, 91) also had their historical figures revised in
#    the same edit (not just shifted) - rewrite D:M for those rows fully.
foreach ($r in $rowOverrides.Keys) {
    $vals = $rowOverrides[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item([int]$r, 4 + $i).Value = $vals[$i]
    }
}
